$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the old D column (min_date), shifting
# min_date/max_date to H/I and mean..lower_quartile to J..O.
$ws.Range("D1:G1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("C1").Value = "n_ascert"
$ws.Range("D1").Value = "n_daily_cases"
$ws.Range("E1").Value = "NAs"
$ws.Range("F1").Value = "n_zeros"
$ws.Range("G1").Value = "n_1to9"

# --- New data columns (n_daily_cases, NAs, n_zeros, n_1to9) ---
# Row 2: Africa
$ws.Range("D2").Value = 44688
$ws.Range("E2").Value = 7197
$ws.Range("F2").Value = 3441
$ws.Range("G2").Value = 12820

# Row 3: Americas
$ws.Range("D3").Value = 27192
$ws.Range("E3").Value = 3939
$ws.Range("F3").Value = 973
$ws.Range("G3").Value = 2211

# Row 4: Asia
$ws.Range("D4").Value = 40831
$ws.Range("E4").Value = 5289
$ws.Range("F4").Value = 2438
$ws.Range("G4").Value = 3590

# Row 5: Europe
$ws.Range("D5").Value = 40609
$ws.Range("E5").Value = 5511
$ws.Range("F5").Value = 751
$ws.Range("G5").Value = 1788

# Row 6: Oceania
$ws.Range("D6").Value = 5049
$ws.Range("E6").Value = 716
$ws.Range("F6").Value = 1240
$ws.Range("G6").Value = 932

# --- Updated min_date / max_date (now uniform across regions) ---
$ws.Range("H2:H6").Value = 43865
$ws.Range("I2:I6").Value = 45017

Write-Host "done"
